# Work Time Table.xlsx update
# - Translate / update the shared strings used in column B (Danish -> English,
#   and correct the 5th entry's duration)
# - Update the selected cell in the sheet view from B5 to B6 (and let the
#   scroll position / topLeftCell reset accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text
$ws.Range("B1").Value = "Amount of work hours"

# Row values - set in the order that matches the target shared-string table
# ordering (B1, B2, B5, B4, B3)
$ws.Range("B2").Value = "4 hours"
$ws.Range("B5").Value = "7 hours 30 minuts"
$ws.Range("B4").Value = "2 hours 30 minuts"
$ws.Range("B3").Value = "4 hours 33 minuts"

# Move the active selection from B5 down to B6
$ws.Range("B6").Select()
